# Weekly update: a new price record for "Comercializadora del Agro de
# Limarí - Haba" is inserted at row 62, shifting the existing rows
# 62-130 down to 63-131 (dimension grows from R130 to R131).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new, blank row above the current row 62. This pushes the
# former row 62 (and everything below it) down by one row, exactly as
# Excel's own "Insert Row" command would.
$ws.Rows.Item(62).Insert()

# Populate the newly inserted row 62 with the new weekly record.
$ws.Cells.Item(62, 1).Value  = 2
$ws.Cells.Item(62, 2).Value  = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(62, 3).Value  = "Coquimbo"
$ws.Cells.Item(62, 4).Value  = 45280
$ws.Cells.Item(62, 5).Value  = 4
$ws.Cells.Item(62, 6).Value  = 100112026
$ws.Cells.Item(62, 7).Value  = "Haba"
$ws.Cells.Item(62, 8).Value  = "Sin especificar"
$ws.Cells.Item(62, 9).Value  = "Primera"
$ws.Cells.Item(62, 10).Value = 400
$ws.Cells.Item(62, 11).Value = 11000
$ws.Cells.Item(62, 12).Value = 13000
$ws.Cells.Item(62, 13).Value = 12000
$ws.Cells.Item(62, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(62, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(62, 16).Value = 480
$ws.Cells.Item(62, 17).Value = 25
$ws.Cells.Item(62, 18).Value = "Hortaliza"
